$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(44977, 45005, 45036, 45066, 45097, 45127, 45158, 45189, 45219, 45250, 45280, 45311)
$newAmount = 8750

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
    $ws.Cells.Item($row, 3).Value = $newAmount
}
